# The commit swaps the two theme parts of this deck:
#   ppt/theme/theme1.xml  (was "Integral", used by the slide master)  -> becomes "Office Theme"
#   ppt/theme/theme2.xml  (was "Office Theme", used by the notes master) -> becomes "Integral"
#
# The PowerPoint object model lets us rewrite the *colour values* of the
# theme that backs the slide master (Master.ColorScheme / RGBColor.RGB),
# so we push the 12 "Office" theme colours into it here. (The scheme/theme
# display names are read-only in the object model - Design.Name has no
# setter in real PowerPoint either - so only the colour values themselves
# are settable through automation.)

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme

# Office theme colours, in PowerPoint Colors(index) order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB() values below are the standard COM 0xBBGGRR encoding of the
# target sRGB hex (e.g. 44546A -> 0x6A5444).
$cs.Colors(1).RGB  = 0x000000   # dk1      000000
$cs.Colors(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$cs.Colors(3).RGB  = 0x6A5444   # dk2      44546A
$cs.Colors(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$cs.Colors(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$cs.Colors(6).RGB  = 0x317DED   # accent2  ED7D31
$cs.Colors(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$cs.Colors(8).RGB  = 0x00C0FF   # accent4  FFC000
$cs.Colors(9).RGB  = 0xC47244   # accent5  4472C4
$cs.Colors(10).RGB = 0x47AD70   # accent6  70AD47
$cs.Colors(11).RGB = 0xC16305   # hlink    0563C1
$cs.Colors(12).RGB = 0x724F95   # folHlink 954F72
